# "Generate Report for Handoff" -- refresh the localization-status report
# with a new pair of source files (old GUIDs -> new GUIDs), flip the
# zh-cn/de-de rows from "handed back" to "ready for handoff", and drop the
# now-stale "latest handback" hyperlinked-filename cells.

$wb = $excel.ActiveWorkbook

$guid1 = "2639c563-4a9d-40ff-a374-978388efcecf"
$guid2 = "ffffc6449859-46b8-4ccd-8643-633cf92b4108"
$xlfHash = "1240c37d6784f132086865b513d001413ae18443"

$xlf1ZhCn = "$guid1.$xlfHash.zh-cn.xlf"
$xlf1DeDe = "$guid1.$xlfHash.de-de.xlf"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$guid1.md"
$wsOverview.Range("C2").Value = ".md"
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-12 03:23:38"

$wsOverview.Range("A3").Value = "$guid2.md"
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-12 03:23:38"

# Update hyperlink display text in place (keep same rIds / target URLs)
$ovSnapshot = @()
foreach ($h in $wsOverview.Hyperlinks) {
    $ovSnapshot += ,@($h.Range.Address(), $h)
}
foreach ($pair in $ovSnapshot) {
    $addr = $pair[0]
    $h = $pair[1]
    if ($addr -eq "`$B`$2") { $h.TextToDisplay = "e2e\$guid1.md" }
    if ($addr -eq "`$B`$3") { $h.TextToDisplay = "e2e\$guid2.md" }
}

# Keep the plain cell text in sync with the hyperlink display text
$wsOverview.Range("B2").Value = "e2e\$guid1.md"
$wsOverview.Range("B3").Value = "e2e\$guid2.md"

$wsOverview.Columns.AutoFit()

# ---------------------------------------------------------------------
# Helper applied to both the "zh-cn" and "de-de" sheets: they share the
# same 16-column layout and the same per-row edits (only the handoff
# xliff filename + datetime text differ between the two languages).
# ---------------------------------------------------------------------
function Update-LangSheet($ws, $xlfName, $handoffDateTime) {
    # Snapshot hyperlinks before mutating anything (Range/Value writes
    # don't reindex the Hyperlinks collection, but Delete() does).
    $snap = @()
    foreach ($h in $ws.Hyperlinks) {
        $snap += ,@($h.Range.Address(), $h)
    }

    # Row 2
    $ws.Range("A2").Value = $guid1 + ".md"
    $ws.Range("C2").Value = "Ready for handoff"
    $ws.Range("G2").Value = $xlfName
    $ws.Range("H2").Value = $handoffDateTime
    $ws.Range("K2").Value = "0001-01-01 00:00:00"

    # Row 3
    $ws.Range("A3").Value = $guid2 + ".md"
    $ws.Range("C3").Value = "Ready for handoff"
    $ws.Range("F3").Value = "True"
    $ws.Range("G3").Value = $xlfName
    $ws.Range("H3").Value = $handoffDateTime
    $ws.Range("K3").Value = "0001-01-01 00:00:00"

    # Refresh the A2/A3 hyperlink display text (still point at the same
    # target rIds / URLs -- only the shown filename changes).
    foreach ($pair in $snap) {
        $addr = $pair[0]
        $h = $pair[1]
        if ($addr -eq "`$A`$2") { $h.TextToDisplay = "$guid1.md" }
        if ($addr -eq "`$A`$3") { $h.TextToDisplay = "$guid2.md" }
    }

    # The "Latest Handback File" (I) / "Latest Handback DateTime" (J)
    # cells no longer apply -- clear them and drop their hyperlinks.
    # Delete hyperlinks back-to-front so earlier indices stay valid.
    for ($i = $snap.Count - 1; $i -ge 0; $i--) {
        $addr = $snap[$i][0]
        $h = $snap[$i][1]
        if ($addr -eq "`$I`$2" -or $addr -eq "`$I`$3") {
            $h.Delete()
        }
    }

    $ws.Range("I2").Style = "Normal"
    $ws.Range("I2").Value = ""
    $ws.Range("J2").Style = "Normal"
    $ws.Range("J2").Value = ""

    $ws.Range("I3").Style = "Normal"
    $ws.Range("I3").Value = ""
    $ws.Range("J3").Style = "Normal"
    $ws.Range("J3").Value = ""

    $ws.Columns.AutoFit()
}

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
Update-LangSheet $wsZhCn $xlf1ZhCn "2016-08-12 03:23:33"

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
Update-LangSheet $wsDeDe $xlf1DeDe "2016-08-12 03:23:38"
